# Update countries & provincias Spain
# Refreshes the COVID snapshot: timestamp, several countries' case counters,
# and three country-name/data swaps caused by the refreshed sort order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 08:54"

# --- Ucrania (row 28) - updated counters only --------------------------
$ws.Range("B28").Value = 244734
$ws.Range("C28").Value = 5397
$ws.Range("D28").Value = 108233
$ws.Range("E28").Value = 131811
$ws.Range("G28").Value = 93
$ws.Range("H28").Value = 4690

# --- Afganistan / Irlanda swap (rows 73-74) -----------------------------
# Afganistan moves up to row 73 with refreshed counters; Irlanda drops to
# row 74 keeping its previous (row 73) counters.
$ws.Range("A73").Value = "Afganistan"
$ws.Range("B73").Value = 39616
$ws.Range("C73").Value = 68
$ws.Range("D73").Value = 33058
$ws.Range("E73").Value = 5088
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 1470

$ws.Range("A74").Value = "Irlanda"
$ws.Range("B74").Value = 39584
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 23364
$ws.Range("E74").Value = 14404
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 1816

# --- El Salvador (row 79) - updated counters only -----------------------
$ws.Range("E79").Value = 4217
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 877

# --- Santa Lucia / Nueva Caledonia swap (rows 207-208) -------------------
# Counters are identical between the two, only the names trade places.
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Nueva Caledonia"

# --- Islas Malvinas / Montserrat swap (rows 215-216) ---------------------
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 1
